$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), cloning the header style used by the
# neighboring "sum" header (G1) so the new column matches the existing
# header formatting (bold, bordered, centered).
$g1 = $ws.Range("G1")
$h1 = $ws.Range("H1")
$g1.Copy($h1)
$h1.Value = "Save"

# New "Save" data value for the single data row (H2).
$ws.Range("H2").Value = 0
